# Started processing factor data into modeling data:
# Insert a new "Include" boolean column (D) into Sheet1, shifting the
# existing Country Currency / Underlying-or-Proxy-Index / REER / Debt-GDP /
# GDP Growth / PMI columns one place to the right (D:I -> E:J).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert the new column before the old column D ("Country Currency").
$ws.Columns.Item(4).Insert()

# Header for the new column.
$ws.Range("D1").Value = "Include"

# Give the new column an explicit (non-autofit) width, matching column C.
$ws.Columns.Item(4).ColumnWidth = $ws.Columns.Item(3).ColumnWidth

# Populate the "Include" flags for each data row (2-48).
$ws.Range("D2:D9").Value = $true
$ws.Range("D10:D13").Value = $false
$ws.Range("D14:D17").Value = $true
$ws.Range("D18").Value = $false
$ws.Range("D19").Value = $true
$ws.Range("D20:D21").Value = $false
$ws.Range("D22:D28").Value = $true
$ws.Range("D29:D30").Value = $false
$ws.Range("D31:D33").Value = $true
$ws.Range("D34:D35").Value = $false
$ws.Range("D36:D45").Value = $true
$ws.Range("D46").Value = $false
$ws.Range("D47:D48").Value = $true

# Collapse the old selection back onto A1 (closest achievable approximation
# of the author's saved view state).
$ws.Range("A1").Select()
